# updateXlsxWithCardsUrl: create a new sheet named with today's date and
# write the cardmarket URLs (with the search filters appended) there, then
# blank out the URL column on the original sheet.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

$sheetName = "27_01_2024"

$newSheet = $wb.Worksheets.Add($null, $ws1)
$newSheet.Name = $sheetName

$lastRow = $ws1.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $a = $ws1.Cells.Item($r, 1).Value2
    $b = $ws1.Cells.Item($r, 2).Value2
    $c = $ws1.Cells.Item($r, 3).Value2
    $d = $ws1.Cells.Item($r, 4).Value2
    $e = $ws1.Cells.Item($r, 5).Value2

    $newSheet.Cells.Item($r, 1).Value = $a
    $newSheet.Cells.Item($r, 2).Value = $b
    $newSheet.Cells.Item($r, 3).Value = $c
    $newSheet.Cells.Item($r, 4).Value = $d

    if ($r -eq 1) {
        $newSheet.Cells.Item($r, 5).Value = ""
    } else {
        if ($d -eq "FR") {
            $langCode = "2"
        } else {
            $langCode = "7"
        }
        $url = "$e" + "?language=" + $langCode + "&minCondition=2&isSigned=N&isPlayset=N&isAltered=N"
        $newSheet.Cells.Item($r, 5).Value = $url
    }

    $ws1.Cells.Item($r, 5).Value = ""
}

$ws1.Select()
